# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-24 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 3
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 2
    21 = 1
    22 = 1
    23 = 0
    24 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
